# This workbook holds a weekly price table. The commit adds one new
# weekly batch of 2 records ("Primera"/Camote and "1a (guarda)"/Paine,
# both sourced from "Perú" on 2021-11-24) at the top of the data block
# (rows 328-329), pushing the rest of the existing records down by two
# rows (the table grows from A1:R424 to A1:R426).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 328; Excel shifts
# all rows from 328 downward (previously 328-424) down to 330-426 and
# copies the row-328 formatting (incl. the date NumberFormat on column D)
# onto the freshly inserted rows.
$ws.Rows("328:329").Insert()

# --- New row 328 ---
$ws.Range("A328").Value = 10
$ws.Range("B328").Value = "Vega Modelo de Temuco"
$ws.Range("C328").Value = "La Araucanía"
$ws.Range("D328").Value = 44524
$ws.Range("E328").Value = 9
$ws.Range("F328").Value = 100112045
$ws.Range("G328").Value = "Zapallo"
$ws.Range("H328").Value = "Camote"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 400
$ws.Range("K328").Value = 800
$ws.Range("L328").Value = 800
$ws.Range("M328").Value = 800
$ws.Range("N328").Value = "$/kilo (volumen en unidades)"
$ws.Range("O328").Value = "Perú"
$ws.Range("P328").Value = 800
$ws.Range("Q328").Value = 1
$ws.Range("R328").Value = "Hortaliza"

# --- New row 329 ---
$ws.Range("A329").Value = 10
$ws.Range("B329").Value = "Vega Modelo de Temuco"
$ws.Range("C329").Value = "La Araucanía"
$ws.Range("D329").Value = 44524
$ws.Range("E329").Value = 9
$ws.Range("F329").Value = 100112045
$ws.Range("G329").Value = "Zapallo"
$ws.Range("H329").Value = "Paine"
$ws.Range("I329").Value = "1a (guarda)"
$ws.Range("J329").Value = 900
$ws.Range("K329").Value = 250
$ws.Range("L329").Value = 300
$ws.Range("M329").Value = 272
$ws.Range("N329").Value = "$/kilo (volumen en unidades)"
$ws.Range("O329").Value = "Perú"
$ws.Range("P329").Value = 272
$ws.Range("Q329").Value = 1
$ws.Range("R329").Value = "Hortaliza"
